# Generate Report for Handoff
# For the four "Ready for handoff" markdown files (rows 4-7) on both the
# zh-cn and de-de localization-status sheets:
#   - bump Priority (column E) from "low" to "ht"
#   - refresh Latest Handoff Datetime (column H) to the new handoff timestamp

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4:E7").Value = "ht"
$wsZhCn.Range("H4:H7").Value = "2016-09-06 02:35:09"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4:E7").Value = "ht"
$wsDeDe.Range("H4:H7").Value = "2016-09-06 02:35:15"

# The Overview sheet's "Latest HO Xliff Generate Date" column mirrors the
# de-de handoff timestamp for these rows, so refresh it too.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4:G7").Value = "2016-09-06 02:35:15"
